$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text shared strings) ---
# A8 rich text ends "...Number  23" -> "...24" (weekly edition number bump)
$a8 = $ws.Range("A8").Value2
$idx = $a8.LastIndexOf("23")
$ws.Range("A8").Characters($idx + 1, 2).Text = "24"

# C9 rich text "Report Covering the Week  6/5/2023  Through  6/11/2023" -> following week
$c9 = $ws.Range("C9").Value2
$idxEnd = $c9.IndexOf("6/11/2023")
$ws.Range("C9").Characters($idxEnd + 1, 9).Text = "6/18/2023"
$idxStart = $c9.IndexOf("6/5/2023")
$ws.Range("C9").Characters($idxStart + 1, 8).Text = "6/12/2023"

# --- Data table updates (rows 15-30): new week of crime-complaint figures ---

# Row 15
$ws.Range("C17").Copy($ws.Range("C15"))
$ws.Range("C15").Value = 1
$ws.Range("C17").Copy($ws.Range("D15"))
$ws.Range("D15").Value = 2
$ws.Range("E17").Copy($ws.Range("E15"))
$ws.Range("E15").Value = -50
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("I15").Value = 18
$ws.Range("J15").Value = 8
$ws.Range("K15").Value = 125
$ws.Range("L15").Value = 260
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = 350

# Row 16
$ws.Range("C16").Value = 11
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 266.666666666667
$ws.Range("F16").Value = 29
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = 70.588235294117
$ws.Range("I16").Value = 129
$ws.Range("J16").Value = 90
$ws.Range("K16").Value = 43.333333333333
$ws.Range("L16").Value = 138.888888888889
$ws.Range("M16").Value = 69.736842105263
$ws.Range("N16").Value = -75.381679389313

# Row 17
$ws.Range("F17").Value = 18
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = -18.181818181818
$ws.Range("I17").Value = 103
$ws.Range("J17").Value = 107
$ws.Range("K17").Value = -3.738317757009
$ws.Range("L17").Value = 19.767441860465
$ws.Range("M17").Value = 66.129032258064
$ws.Range("N17").Value = -21.969696969697

# Row 18
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 40
$ws.Range("I18").Value = 87
$ws.Range("J18").Value = 83
$ws.Range("K18").Value = 4.819277108433
$ws.Range("L18").Value = 12.987012987013
$ws.Range("M18").Value = -26.271186440678
$ws.Range("N18").Value = -87.053571428571

# Row 19
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -8.333333333333
$ws.Range("F19").Value = 56
$ws.Range("G19").Value = 44
$ws.Range("H19").Value = 27.272727272727
$ws.Range("I19").Value = 335
$ws.Range("J19").Value = 298
$ws.Range("K19").Value = 12.416107382550
$ws.Range("L19").Value = 77.248677248677
$ws.Range("M19").Value = 64.215686274509
$ws.Range("N19").Value = -15.404040404040

# Row 20
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 150
$ws.Range("F20").Value = 16
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = 23.076923076923
$ws.Range("I20").Value = 100
$ws.Range("J20").Value = 83
$ws.Range("K20").Value = 20.481927710843
$ws.Range("L20").Value = 25
$ws.Range("M20").Value = -7.407407407407
$ws.Range("N20").Value = -88.713318284424

# Row 21
$ws.Range("C21").Value = 33
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = 26.923076923076
$ws.Range("F21").Value = 134
$ws.Range("G21").Value = 108
$ws.Range("H21").Value = 24.074074074074
$ws.Range("I21").Value = 772
$ws.Range("J21").Value = 669
$ws.Range("K21").Value = 15.396113602391
$ws.Range("L21").Value = 56.910569105691
$ws.Range("M21").Value = 33.333333333333
$ws.Range("N21").Value = -70.545593285005

# Row 22
$ws.Range("C22").Value = 3
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 200
$ws.Range("F22").Value = 7
$ws.Range("G22").Value = 10
$ws.Range("H22").Value = -30
$ws.Range("I22").Value = 43
$ws.Range("J22").Value = 38
$ws.Range("K22").Value = 13.157894736842
$ws.Range("L22").Value = 186.666666666667
$ws.Range("M22").Value = 138.888888888889

# Row 24
$ws.Range("C24").Value = 36
$ws.Range("D24").Value = 38
$ws.Range("E24").Value = -5.263157894736
$ws.Range("F24").Value = 152
$ws.Range("G24").Value = 118
$ws.Range("H24").Value = 28.813559322033
$ws.Range("I24").Value = 794
$ws.Range("J24").Value = 611
$ws.Range("K24").Value = 29.950900163666
$ws.Range("L24").Value = 47.037037037037
$ws.Range("M24").Value = 86.384976525821

# Row 25
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 33
$ws.Range("G25").Value = 46
$ws.Range("H25").Value = -28.260869565217
$ws.Range("I25").Value = 239
$ws.Range("J25").Value = 265
$ws.Range("K25").Value = -9.811320754716
$ws.Range("L25").Value = 24.479166666666
$ws.Range("M25").Value = 2.575107296137

# Row 26
$ws.Range("C17").Copy($ws.Range("C26"))
$ws.Range("C26").Value = 2
$ws.Range("C17").Copy($ws.Range("D26"))
$ws.Range("D26").Value = 2
$ws.Range("E17").Copy($ws.Range("E26"))
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 3
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 23
$ws.Range("J26").Value = 11
$ws.Range("K26").Value = 109.090909090909
$ws.Range("L26").Value = 130

# Row 27
$ws.Range("G27").Value = 9
$ws.Range("H27").Value = -11.111111111111
$ws.Range("I27").Value = 49
$ws.Range("J27").Value = 43
$ws.Range("K27").Value = 13.953488372093
$ws.Range("L27").Value = 75

# Row 28
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = -50
$ws.Range("J28").Value = 3
$ws.Range("K28").Value = 0

# Row 29
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = -50
$ws.Range("J29").Value = 3
$ws.Range("K29").Value = 0

# Row 30
$ws.Range("C23").Copy($ws.Range("D30"))
$ws.Range("E23").Copy($ws.Range("E30"))
$ws.Range("G30").Value = 1
